$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update inputs per "Spinney and Matt updates and LHS implemented."
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 15
$ws.Range("M4").Value = 17
